$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

$ws1.Name = "Sheet1 - Text"
$ws2.Name = "Sheet2 - Numbers"
$ws3.Name = "Sheet3 - Formulas"

# --- Sheet1 - Text: text data ---
$ws1.Range("A1").Value = "This is cell A1 in Sheet 1"
$ws1.Range("G5").Value = "This is cell G5"
[void]$ws1.Range("G6").Select()

# --- Sheet2 - Numbers: numbers in D1:D30, percentages in K1:K30 ---
for ($i = 1; $i -le 30; $i++) {
    $ws2.Cells.Item($i, 4).Value = $i
}

$ws2.Range("G5").Value = "This is cell G5"

for ($i = 1; $i -le 30; $i++) {
    $ws2.Cells.Item($i, 11).Value = $i / 100
}
$ws2.Range("K1:K30").Style = "Percent"

[void]$ws2.Range("L2").Select()

# --- Sheet3 - Formulas: formula referencing Sheet2 - Numbers ---
$ws3.Range("D2").Formula = "='Sheet2 - Numbers'!D5"
[void]$ws3.Range("D3").Select()

# --- Active sheet / tab selection ---
[void]$ws3.Activate()

$wb.Save()
